$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to remain plain text so values like "1.00" or "3.90"
# are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '88.412.42'
$ws.Range('E2').Value = '  -4.33%  '

$ws.Range('D3').Value = '3.121.57'
$ws.Range('E3').Value = '  -5.07%  '

$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.48%  '

$ws.Range('D5').Value = '213.09'
$ws.Range('E5').Value = '  -0.32%  '

$ws.Range('D6').Value = '632.39'
$ws.Range('E6').Value = '  +0.77%  '

$ws.Range('D7').Value = '0.392'
$ws.Range('E7').Value = '  -5.88%  '

$ws.Range('D8').Value = '0.730'
$ws.Range('E8').Value = '  +3.24%  '

$ws.Range('E9').Value = '  +0.02%  '

$ws.Range('D10').Value = '3.119.07'
$ws.Range('E10').Value = '  -5.19%  '

$ws.Range('D11').Value = '0.553'
$ws.Range('E11').Value = '  -5.40%  '

$ws.Range('D12').Value = '0.179'
$ws.Range('E12').Value = '  -0.65%  '

$ws.Range('D13').Value = '0.0000249'
$ws.Range('E13').Value = '  -5.70%  '

$ws.Range('D14').Value = '5.28'
$ws.Range('E14').Value = '  -1.58%  '

$ws.Range('D15').Value = '88.240.22'
$ws.Range('E15').Value = '  -4.31%  '

$ws.Range('D16').Value = '3.687.74'
$ws.Range('E16').Value = '  -5.83%  '

$ws.Range('D17').Value = '32.19'
$ws.Range('E17').Value = '  -6.24%  '

$ws.Range('D18').Value = '3.117.05'
$ws.Range('E18').Value = '  -5.87%  '

$ws.Range('D19').Value = '3.33'
$ws.Range('E19').Value = '  +1.84%  '

$ws.Range('D20').Value = '0.0000217'
$ws.Range('E20').Value = '  +16.83%  '

$ws.Range('D21').Value = '13.16'
$ws.Range('E21').Value = '  -6.37%  '

$ws.Range('D22').Value = '425.33'
$ws.Range('E22').Value = '  -3.27%  '

$ws.Range('D23').Value = '8.33'
$ws.Range('E23').Value = '  -6.04%  '

$ws.Range('D24').Value = '4.87'
$ws.Range('E24').Value = '  -7.72%  '

$ws.Range('D25').Value = '5.35'
$ws.Range('E25').Value = '  +0.07%  '

$ws.Range('D26').Value = '11.45'
$ws.Range('E26').Value = '  -6.97%  '

$ws.Range('D27').Value = '79.54'
$ws.Range('E27').Value = '  +4.27%  '

$ws.Range('D28').Value = '3.270.70'
$ws.Range('E28').Value = '  -7.32%  '

$ws.Range('E29').Value = '  +0.16%  '

$ws.Range('D30').Value = '0.159'
$ws.Range('E30').Value = '  -11.45%  '

$ws.Range('B31').Value = 'Binance-PegBSC-USD'
$ws.Range('C31').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D31').Value = '0.993'
$ws.Range('E31').Value = '  -0.61%  '

$ws.Range('B32').Value = 'dogwifhat'
$ws.Range('C32').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D32').Value = '3.90'
$ws.Range('E32').Value = '  +7.38%  '

$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').Value = '8.15'
$ws.Range('E33').Value = '  -8.00%  '

$ws.Range('B34').Value = 'Bittensor'
$ws.Range('C34').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D34').Value = '511.88'
$ws.Range('E34').Value = '  -7.57%  '

$ws.Range('B35').Value = 'RenderToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D35').Value = '6.93'
$ws.Range('E35').Value = '  -3.50%  '

$ws.Range('B36').Value = 'Fetch.AI'
$ws.Range('C36').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D36').Value = '1.28'
$ws.Range('E36').Value = '  -1.30%  '

$ws.Range('B37').Value = 'PancakeSwap'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D37').Value = '1.83'
$ws.Range('E37').Value = '  -4.72%  '

$ws.Range('D38').Value = '21.79'
$ws.Range('E38').Value = '  -4.13%  '

$ws.Range('B39').Value = 'WhiteBITCoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D39').Value = '22.19'
$ws.Range('E39').Value = '  -1.21%  '

$ws.Range('B40').Value = 'FirstDigitalUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  +0.04%  '

$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = '0.125'
$ws.Range('E41').Value = '  -5.15%  '

$ws.Range('B42').Value = 'USDe'
$ws.Range('C42').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  -0.11%  '

$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').Value = '1.85'
$ws.Range('E43').Value = '  -6.73%  '

$ws.Range('B44').Value = 'PolygonEcosystemToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D44').Value = '0.361'
$ws.Range('E44').Value = '  -8.44%  '

$ws.Range('B45').Value = 'Monero'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D45').Value = '143.71'
$ws.Range('E45').Value = '  -4.56%  '

$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D46').Value = '43.61'
$ws.Range('E46').Value = '  -0.22%  '

$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D47').Value = '0.126'
$ws.Range('E47').Value = '  -3.04%  '

$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').Value = '164.49'
$ws.Range('E48').Value = '  -9.18%  '

$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').Value = '0.719'
$ws.Range('E49').Value = '  -1.75%  '

$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').Value = '24.27'
$ws.Range('E50').Value = '  -3.09%  '

$ws.Range('B51').Value = 'ImmutableX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D51').Value = '1.18'
$ws.Range('E51').Value = '  -8.08%  '
